# Apply the StructureDefinition-employee-department.xlsx update:
#  - Metadata sheet: bump Version, Date, add Publisher/Jurisdiction info,
#    drop the duplicated "Contact" row (net: 21 -> 20 rows).
#  - Elements sheet: give the root Extension row its own Short/Definition
#    text instead of the generic "Extension"/"An Extension" boilerplate.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Metadata ----
$meta = $wb.Worksheets.Item(1)

# Version 5.0.0 -> 6.0.0
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date -> new publish date
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction"
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row; remove it
# entirely (shifts Description..Context up by one row to match the target layout).
$meta.Rows.Item(11).Delete()

# ---- Sheet 2: Elements ----
$elements = $wb.Worksheets.Item(2)

# The root "Extension" element row (row 2) now carries its own Short/Definition
# text (reusing the Title/Description wording) instead of the generic
# "Extension" / "An Extension" placeholders.
$elements.Cells.Item(2, 11).Value = "Employee Department"
$elements.Cells.Item(2, 12).Value = "Code for the department of the employee"
